# remove expbottle, fix some small bugs
#
# The "MainIcon" sheet has a table (表1) listing UI icons. The row for the
# "经验" / "打开经验瓶" (expbottle) icon (Id=40, currently worksheet row 20)
# is being removed entirely. Deleting the whole row shifts the rows below
# it up by one (so the old row 21 becomes row 20, old row 22 becomes row
# 21), shrinks the table/dimension from L22 to L21, and the now-unused
# shared strings ("经验", "打开经验瓶", "SideButton15") drop out of
# sharedStrings.xml automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MainIcon")

# Delete the entire row that holds the expbottle entry (Id = 40).
$ws.Rows.Item(20).EntireRow.Delete()

# Match the author's final cursor position recorded in the saved file.
$ws.Range("E11").Select()
